# Auto-generated edit script: updates Price (D) and Volume(1h) (E) columns
# of the cryptos sheet with refreshed values, preserving cell text formatting.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.017.25"
$ws.Range("E2").Value = "  -2.42%  "
$ws.Range("D3").Value = "3.519.39"
$ws.Range("E3").Value = "  -3.17%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "586.70"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.46%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "169.91"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.48%  "
$ws.Range("E7").Value = "  -0.89%  "
$ws.Range("D8").Value = "3.516.42"
$ws.Range("E8").Value = "  -3.02%  "
$ws.Range("E9").Value = "  -0.02%  "
$ws.Range("E10").Value = "  -4.20%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.80"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.31%  "
$ws.Range("E12").Value = "  -5.27%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "47.47"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.19%  "
$ws.Range("E14").Value = "  -3.23%  "
$ws.Range("D15").Value = "4.088.78"
$ws.Range("E15").Value = "  -3.15%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "8.42"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -6.32%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "612.33"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -9.21%  "
$ws.Range("D18").Value = "69.108.22"
$ws.Range("E18").Value = "  -2.38%  "
$ws.Range("D19").Value = "3.517.76"
$ws.Range("E19").Value = "  -3.22%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.120"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.66%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.38"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.31%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.07"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.75%  "
$ws.Range("E23").Value = "  -6.15%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "15.75"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -8.40%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "96.50"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.40%  "
$ws.Range("E26").Value = "  -2.41%  "
$ws.Range("E27").Value = "  +0.00%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.61"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -6.42%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.21"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -6.12%  "
$ws.Range("E30").Value = "  -5.78%  "
$ws.Range("E31").Value = "  -7.10%  "
$ws.Range("E32").Value = "  -5.14%  "
$ws.Range("E33").Value = "  -4.78%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.91"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -8.45%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "614.75"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +6.74%  "
$ws.Range("E36").Value = "  -3.24%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.48"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -12.41%  "
$ws.Range("E38").Value = "  -5.39%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "57.04"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.51%  "
$ws.Range("E40").Value = "  +0.02%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0445"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.51%  "
$ws.Range("D42").Value = "3.399.48"
$ws.Range("E42").Value = "  -4.15%  "
$ws.Range("E43").Value = "  -3.55%  "
$ws.Range("E44").Value = "  -5.55%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "32.80"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.69%  "
$ws.Range("D46").Value = "0.0₃0697"
$ws.Range("E46").Value = "  -5.02%  "
$ws.Range("E47").Value = "  -5.69%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.74"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -7.44%  "
$ws.Range("E49").Value = "  -3.60%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "133.91"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.84%  "
$ws.Range("E51").Value = "  +10.17%  "
